$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New challenge column: "sorting-by-bits" added after the existing M column.
# Header cell N1 should look exactly like the other header cells (same
# style as M1) - copy M1 (value+format) into N1, then overwrite its value.
$ws.Range("M1").Copy($ws.Range("N1"))
$ws.Range("N1").Value = "sorting-by-bits"

# Every contestant's result for the new challenge is recorded as boolean
# FALSE (no style, like the rest of the data cells).
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 14).Value = $false
}
